$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A:A").Insert()

$ws.Range("A1").Value = "senaryoturu"
$ws.Range("A2").Value = "negatif"
$ws.Range("A3").Value = "negatif"
$ws.Range("A4").Value = "pozitif"
$ws.Range("B4").Value = "diamond"

$ws.Range("C4").Select() | Out-Null
